$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 253
$ws1.Range("F6").Value = 559
$ws1.Range("F9").Value = 274
$ws1.Range("F10").Value = 387
$ws1.Range("F12").Value = 699
$ws1.Range("F13").Value = 762
$ws1.Range("F14").Value = 1518
$ws1.Range("F15").Value = 1518
$ws1.Range("F16").Value = 892
$ws1.Range("F17").Value = 30
$ws1.Range("F20").Value = 329
$ws1.Range("F24").Value = 6635
$ws1.Range("F25").Value = 4998
$ws1.Range("F28").Value = 209
$ws1.Range("F29").Value = 203
$ws1.Range("F32").Value = 1286
$ws1.Range("F33").Value = 196
$ws1.Range("F34").Value = 252
$ws1.Range("F35").Value = 617
$ws1.Range("F38").Value = 250
$ws1.Range("F40").Value = 149

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2464
$ws3.Range("F4").Value = 200
$ws3.Range("F5").Value = 61

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 253
$ws4.Range("F6").Value = 200
$ws4.Range("F7").Value = 61
$ws4.Range("F9").Value = 559
$ws4.Range("F12").Value = 274
$ws4.Range("F14").Value = 387
$ws4.Range("F16").Value = 699
$ws4.Range("F17").Value = 762
$ws4.Range("F18").Value = 1518
$ws4.Range("F19").Value = 1518
$ws4.Range("F20").Value = 892
$ws4.Range("F21").Value = 30
$ws4.Range("F24").Value = 329
$ws4.Range("F29").Value = 6635
$ws4.Range("F30").Value = 4998
$ws4.Range("F32").Value = 203
$ws4.Range("F34").Value = 1286
$ws4.Range("F35").Value = 196
$ws4.Range("F36").Value = 252
$ws4.Range("F39").Value = 617
$ws4.Range("F43").Value = 250
$ws4.Range("F44").Value = 149
